$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, pushing the existing rows 232-247 down to 233-248
$ws.Rows.Item(232).Insert()

# Fill in the new row 232 with the new weekly record
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44585
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 100114014
$ws.Range("G232").Value = "Betarraga"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 4000
$ws.Range("K232").Value = 500
$ws.Range("L232").Value = 500
$ws.Range("M232").Value = 500
$ws.Range("N232").Value = "`$/paquete 5 unidades"
$ws.Range("O232").Value = "Región del Maule"
$ws.Range("P232").Value = 100
$ws.Range("Q232").Value = 5
$ws.Range("R232").Value = "Hortaliza"
